$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q4" sheet, positioned right after "总计" and before
#    "2022-Q1". The easiest way to get identical structure/formatting to the
#    existing quarter sheets is to duplicate "2022-Q1" and then overwrite the
#    numbers that differ.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1)                      # places the copy immediately before $q1
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The D:G columns hold their numbers as plain text (matching the source
# sheet), so force text formatting before writing the new values; otherwise
# Excel would silently reinterpret "1.06" etc. as a number.
$q4.Range("D2:G4").NumberFormat = "@"

$q4.Range("D2").Value = "1.06"
$q4.Range("E2").Value = "82.28"
$q4.Range("F2").Value = "1.92"
$q4.Range("G2").Value = "0.0204"
$q4.Range("H2").Value = 10

$q4.Range("D3").Value = "1.06"
$q4.Range("E3").Value = "82.28"
$q4.Range("F3").Value = "1.92"
$q4.Range("G3").Value = "0.0204"
$q4.Range("H3").Value = 10

$q4.Range("D4").Value = "1.06"
$q4.Range("E4").Value = "82.28"
$q4.Range("F4").Value = "1.92"
$q4.Range("G4").Value = "0.0204"
$q4.Range("H4").Value = 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: push the existing two data rows down
#    by one and insert a new row for "2022-Q4" on top, keeping the per-cell
#    formatting intact (column A carries a distinct style).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift row 3 ("2021-Q4") down to row 4.
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)   # xlPasteFormats
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.38

# Shift row 2 ("2022-Q1") down to row 3.
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)   # xlPasteFormats
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.02

# Write the new "2022-Q4" row into row 2 (format is already in place).
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.06
